$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The series (daily new-positives / 7-day rolling sum / rate per 100k) continues
# past the previous last day (2021-09-20, row 385) through 2021-12-08 (row 464).
# First extend the row formatting (column A date style s="2", thin borders, bold
# header font, center/top alignment) down to the new rows by copying the last
# existing data row.
$ws.Range("A385:D385").Copy()
$ws.Range("A386:D464").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Daily data: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44460, 0, 3, 250.208507089241),
    @(44461, 0, 3, 250.208507089241),
    @(44462, 0, 3, 250.208507089241),
    @(44463, 1, 1, 83.40283569641367),
    @(44464, 0, 1, 83.40283569641367),
    @(44465, 0, 1, 83.40283569641367),
    @(44466, 1, 2, 166.8056713928273),
    @(44467, 0, 2, 166.8056713928273),
    @(44468, 0, 2, 166.8056713928273),
    @(44469, 1, 3, 250.208507089241),
    @(44470, 2, 4, 333.6113427856547),
    @(44471, 0, 4, 333.6113427856547),
    @(44472, 0, 4, 333.6113427856547),
    @(44473, 1, 4, 333.6113427856547),
    @(44474, 0, 4, 333.6113427856547),
    @(44475, 0, 4, 333.6113427856547),
    @(44476, 0, 3, 250.208507089241),
    @(44477, 0, 1, 83.40283569641367),
    @(44478, 0, 1, 83.40283569641367),
    @(44479, 0, 1, 83.40283569641367),
    @(44480, 0, 0, 0),
    @(44481, 0, 0, 0),
    @(44482, 0, 0, 0),
    @(44483, 0, 0, 0),
    @(44484, 0, 0, 0),
    @(44485, 0, 0, 0),
    @(44486, 0, 0, 0),
    @(44487, 0, 0, 0),
    @(44488, 0, 0, 0),
    @(44489, 0, 0, 0),
    @(44490, 0, 0, 0),
    @(44491, 0, 0, 0),
    @(44492, 0, 0, 0),
    @(44493, 0, 0, 0),
    @(44494, 0, 0, 0),
    @(44495, 0, 0, 0),
    @(44496, 0, 0, 0),
    @(44497, 0, 0, 0),
    @(44498, 0, 0, 0),
    @(44499, 0, 0, 0),
    @(44500, 0, 0, 0),
    @(44501, 0, 0, 0),
    @(44502, 0, 0, 0),
    @(44503, 0, 0, 0),
    @(44504, 0, 0, 0),
    @(44505, 0, 0, 0),
    @(44506, 0, 0, 0),
    @(44507, 0, 0, 0),
    @(44508, 0, 0, 0),
    @(44509, 0, 0, 0),
    @(44510, 0, 0, 0),
    @(44511, 0, 0, 0),
    @(44512, 0, 0, 0),
    @(44513, 0, 0, 0),
    @(44514, 0, 0, 0),
    @(44515, 0, 0, 0),
    @(44516, 0, 0, 0),
    @(44517, 0, 0, 0),
    @(44518, 0, 0, 0),
    @(44519, 0, 0, 0),
    @(44520, 0, 0, 0),
    @(44521, 0, 0, 0),
    @(44522, 0, 0, 0),
    @(44523, 0, 0, 0),
    @(44524, 0, 0, 0),
    @(44525, 0, 0, 0),
    @(44526, 0, 0, 0),
    @(44527, 0, 0, 0),
    @(44528, 0, 0, 0),
    @(44529, 2, 2, 166.8056713928273),
    @(44530, 0, 2, 166.8056713928273),
    @(44531, 0, 2, 166.8056713928273),
    @(44532, 1, 3, 250.208507089241),
    @(44533, 3, 6, 500.4170141784821),
    @(44534, 0, 6, 500.4170141784821),
    @(44535, 3, 9, 750.6255212677231),
    @(44536, 1, 8, 667.2226855713094),
    @(44537, 0, 8, 667.2226855713094),
    @(44538, 1, 9, 750.6255212677231),
)

$startRow = 386
$rowCount = $data.Count
$arr = New-Object "object[,]" $rowCount,4
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}
$destination = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $rowCount - 1, 4))
$destination.Value = $arr

Write-Output ("rows written: " + $rowCount)